$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated b0 (B) and b1 (C) parameter values for rows 267-281 (IDs 385-400)
# as part of the "renew parameters in 15th Jan 2026" update.

$updates = @(
    @{ Row = 267; B = -0.0196; C = 0.0509 },
    @{ Row = 268; B = -0.0055; C = 0.1014 },
    @{ Row = 269; B = 0.0239;  C = 0.1708 },
    @{ Row = 270; B = 0.0418;  C = 0.2056 },
    @{ Row = 271; B = 0.003;   C = 0.433 },
    @{ Row = 272; B = 0.0272;  C = 0.1196 },
    @{ Row = 273; B = 0.0093;  C = 0.1519 },
    @{ Row = 274; B = 0.0036;  C = 0.2025 },
    @{ Row = 275; B = 0.0956;  C = 0.3232 },
    @{ Row = 276; B = 0.0989;  C = 0.1515 },
    @{ Row = 277; B = 0.0046;  C = 0.1715 },
    @{ Row = 278; B = -0.0136; C = 0.0814 },
    @{ Row = 279; B = -0.0123; C = 0.1297 },
    @{ Row = 280; B = 0.0504;  C = 0.3509 },
    @{ Row = 281; B = 0.0133;  C = 0.1265 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}
